$d = $word.ActiveDocument
$newText = "Perioadele campaniei din Gemini: 14-23 februarie, 14-24 martie"

function Replace-Once($pattern) {
    $rng = $d.Content
    $found = $rng.Find.Execute(
        $pattern, $true, $false, $false, $false, $false, $true, 1, $false,
        "", 1)
    if (-not $found) {
        throw "Pattern not found: $pattern"
    }
    $rng.InsertAfter($newText)
    $rng.Font.Reset()
}

# Occurrence 1: leading red space run + "Perioadele campaniei din 2018 pentru " +
# "Perseu" + ": " + "30 octombrie-8 noiembrie și 29 noiembrie-8 decembrie" runs.
Replace-Once(" Perioadele campaniei din 2018 pentru Perseu: 30 octombrie-8 noiembrie și 29 noiembrie-8 decembrie")

# Occurrence 2: single text run + trailing space run.
Replace-Once("Perioadele campaniei din 2018 pentru Perseu: 30 octombrie-8 noiembrie și 29 noiembrie-8 decembrie ")

# Occurrence 3: single text run, no surrounding whitespace runs.
Replace-Once("Perioadele campaniei din 2018 pentru Perseu: 30 octombrie-8 noiembrie și 29 noiembrie-8 decembrie")

# Occurrence 4: single text run, no surrounding whitespace runs.
Replace-Once("Perioadele campaniei din 2018 pentru Perseu: 30 octombrie-8 noiembrie și 29 noiembrie-8 decembrie")

Write-Host "All 4 occurrences replaced."
